# 9th Stab - Cosmetic Changes
# Add two new "watch" date columns (Jun_15, Jun_17) ahead of the existing
# Jun_13 column, shifting the old Jun_13/Jun_10 columns two slots to the
# right, and seed the new columns with the "UN" default rating used
# throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at B:C - this pushes the existing column B (Jun_13
# header + data) to D, and existing column C (Jun_10 header + data) to E.
$ws.Columns("B:C").Insert()

# Header row: set C1 before B1 so the new shared-strings are registered in
# the same order (Jun_15 then Jun_17) as the target workbook.
$ws.Range("C1").Value = "Jun_15"
$ws.Range("B1").Value = "Jun_17"

# Data rows: the new columns default to the same "UN" rating as the rest of
# the sheet for every existing data row (rows 2-27).
$ws.Range("B2:C27").Value = "UN"
